$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-5 (odoo connection block): "Usar" flips from True to False.
# Copy from an existing text "False" cell (C6) so the value lands as shared
# text ("False") rather than being auto-typed to a native boolean.
$ws.Range("C6").Copy($ws.Range("C2"))
$ws.Range("C6").Copy($ws.Range("C3"))
$ws.Range("C6").Copy($ws.Range("C4"))
$ws.Range("C6").Copy($ws.Range("C5"))

# --- New rows 12-15: server parameters for the new Efectores / SIF-SIGEHOS join
$ws.Range("A12").Value = "host"
$ws.Range("B12").Value = "10.22.0.142"

$ws.Range("A13").Value = "user"
$ws.Range("B13").Value = "postgres"

$ws.Range("A14").Value = "password"
$ws.Range("B14").Value = "serveradmin"

$ws.Range("A15").Value = "database"
$ws.Range("B15").Value = "Facoep"

# "Usar" column for the new rows is True -- copy from an existing text
# "True" cell (C11) for the same reason as above.
$ws.Range("C11").Copy($ws.Range("C12"))
$ws.Range("C11").Copy($ws.Range("C13"))
$ws.Range("C11").Copy($ws.Range("C14"))
$ws.Range("C11").Copy($ws.Range("C15"))

# B12 (host) and B14 (password) get a distinct font: Arial 10.
# Apply to both cells as one union range so the style is computed once.
$fontRng = $ws.Range("B12,B14")
$fontRng.Font.Name = "Arial"
$fontRng.Font.Size = 10

# Grow the table to cover the newly-added rows
$ws.ListObjects.Item(1).Resize($ws.Range("A1:C15"))

# Move selection to the first empty row below the table, matching the saved state
$ws.Range("A16").Select()
